# Login.xlsx: "Persoon" -> "Account" rename
#
# The "Persoon" worksheet's first two columns hold a label row (A1) and a
# value row (A2) that spell out the entity name in two forms:
#   A1 = "[Persoon]"   (bracketed / token form)
#   A2 = "Persoon"     (plain form)
# The commit renames that entity to "Account" everywhere it is referenced
# on this sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Persoon")

# Write A2 ("Account") before A1 ("[Account]") so the new shared-string
# table entries land in the same order Excel itself produced them in.
$ws.Range("A2").Value = "Account"
$ws.Range("A1").Value = "[Account]"

# Leave the active cell on A2, matching the saved selection state.
[void]$ws.Range("A2").Select()
